# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the 5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.* file across the Overview,
# zh-cn and de-de sheets, as produced by a fresh handback-status report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
# "Latest HO Xliff Generate Date" for the 5ffd0142... row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-18 06:46:32"

# --- zh-cn sheet -------------------------------------------------------
# "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K)
# for the 5ffd0142... row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-18 06:46:26"
$wsZhCn.Range("K3").Value = "2016-08-18 06:46:53"

# --- de-de sheet ---------------------------------------------------------
# "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K)
# for the 5ffd0142... row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-18 06:46:32"
$wsDeDe.Range("K3").Value = "2016-08-18 06:47:04"
